# Applies the "Last changes in the project" edit:
#  - Adds three new columns T, U, V (1/0 indicator flags) for every data
#    row (1..150), following a repeating 10-row cycle.
#  - Highlights the T:V cells of every row whose row-in-cycle position is 3
#    (rows 3,13,23,...,143) with a solid yellow fill - this creates the new
#    fill/cellXf style entries.
#  - Moves the active selection from G7 to A123.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repeating 10-row cycle of (T,U,V) values, indexed 1..10 (index 0 unused).
$tCycle = @(0,0,0,0,0,0,1,0,1,0,0)
$uCycle = @(0,0,0,1,0,0,0,1,0,1,0)
$vCycle = @(0,1,1,0,1,1,0,0,0,0,1)

$lastRow = 150

for ($r = 1; $r -le $lastRow; $r++) {
    $pos = $r % 10
    if ($pos -eq 0) { $pos = 10 }

    $tCell = $ws.Cells.Item($r, 20)
    $uCell = $ws.Cells.Item($r, 21)
    $vCell = $ws.Cells.Item($r, 22)

    $tCell.Value = $tCycle[$pos]
    $uCell.Value = $uCycle[$pos]
    $vCell.Value = $vCycle[$pos]

    if ($pos -eq 3) {
        $ws.Range($tCell, $vCell).Interior.Color = 65535
    }
}

# Move the selection, matching the new activeCell/sqref in the diff.
[void]$ws.Range("A123").Select()
